$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, -0.1349883557924783, 1.93267810644422, 16.79600653627128, 4.098293124737575, 4.19246071033399, 22),
    @(3, 0.1207088239491281, 1.786104686133971, 11.46089906083781, 3.385394963787507, 3.46679172473507, 21),
    @(4, -0.5555709582211426, 1.050109431406426, 4.446450755259711, 2.108660891480589, 2.087000281350939, 20),
    @(5, 0.07663410097718948, 0.7265600105831413, 1.828123001807706, 1.352080989367022, 1.386898099900304, 19),
    @(6, 0.04069671201648461, 0.7470070225631975, 1.696764068754185, 1.302598966971103, 1.339708950721457, 18),
    @(7, -0.003406094425567046, 0.6048042399408757, 0.6518029153055619, 0.8073431211731242, 0.8321828350431602, 17),
    @(8, 0.08389501441541836, 0.5881369787635709, 0.6090337700840233, 0.7804061571284681, 0.8013291602599671, 16)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
